$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 11666
# Row 64
$ws.Range("H64").Value = 4072.1667
$ws.Range("I64").Value = 4070.5293
$ws.Range("J64").Value = 4100
$ws.Range("K64").Value = 4070.5293
$ws.Range("L64").Value = 4100
$ws.Range("M64").Value = -3822.5293
$ws.Range("N64").Value = -4596
# Row 67
$ws.Range("H67").Value = 4072.1667
$ws.Range("I67").Value = 4070.5293
$ws.Range("J67").Value = 4100
$ws.Range("K67").Value = 4070.5293
$ws.Range("L67").Value = 4100
$ws.Range("M67").Value = -3212.5293
$ws.Range("N67").Value = -5816

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 41.2
$ws.Range("I5").Value = 41.5
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 41.5
$ws.Range("L5").Value = 40
$ws.Range("M5").Value = 70.5
$ws.Range("N5").Value = -264
# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
# Row 45
$ws.Range("H45").Value = 9557.5
$ws.Range("I45").Value = 2116.6667
$ws.Range("K45").Value = 2116.6667
$ws.Range("M45").Value = -1739.6667
# Row 61
$ws.Range("H61").Value = 12211.111
$ws.Range("I61").Value = 10816.667
$ws.Range("J61").Value = 15000
$ws.Range("K61").Value = 10816.667
$ws.Range("L61").Value = 15000
$ws.Range("M61").Value = -10604.667
$ws.Range("N61").Value = -15424
# Row 136
$ws.Range("H136").Value = 12211.111
$ws.Range("I136").Value = 10816.667
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 32450.001
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -29900.001
$ws.Range("N136").Value = -50100

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 41.2
$ws.Range("I4").Value = 41.5
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 41.5
$ws.Range("L4").Value = 40
$ws.Range("M4").Value = 73.5
$ws.Range("N4").Value = -270
# Row 22
$ws.Range("H22").Value = 775
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -127
# Row 86
$ws.Range("H86").Value = 1005
$ws.Range("I86").Value = 1082
$ws.Range("J86").Value = 389
$ws.Range("K86").Value = 1082
$ws.Range("L86").Value = 389
$ws.Range("M86").Value = 41
$ws.Range("N86").Value = -2635
# Row 89
$ws.Range("H89").Value = 1005
$ws.Range("I89").Value = 1082
$ws.Range("J89").Value = 389
$ws.Range("K89").Value = 5410
$ws.Range("L89").Value = 1945
$ws.Range("M89").Value = 206
$ws.Range("N89").Value = -13177
# Row 130
$ws.Range("H130").Value = 40375
$ws.Range("J130").Value = 40375
$ws.Range("L130").Value = 40375
$ws.Range("N130").Value = -50415

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 2033.3334
$ws.Range("J3").Value = 2033.3334
$ws.Range("L3").Value = 2033.3334
$ws.Range("N3").Value = -2259.3334
# Row 69
$ws.Range("H69").Value = 18045.5
$ws.Range("I69").Value = 6091
$ws.Range("K69").Value = 6091
$ws.Range("M69").Value = -5342
# Row 72
$ws.Range("H72").Value = 18045.5
$ws.Range("I72").Value = 6091
$ws.Range("K72").Value = 18273
$ws.Range("M72").Value = -14529
# Row 86
$ws.Range("H86").Value = 3750
$ws.Range("I86").Value = 3750
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3750
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2627
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 3750
$ws.Range("I89").Value = 3750
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 18750
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -13134
$ws.Range("N89").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 465.77777
$ws.Range("I34").Value = 132.5
$ws.Range("J34").Value = 1132.3334
$ws.Range("K34").Value = 397.5
$ws.Range("L34").Value = 3397.0002
$ws.Range("M34").Value = -313.5
$ws.Range("N34").Value = -3565.0002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 18.714285
$ws.Range("J2").Value = 19.714285
$ws.Range("L2").Value = 19.714285
$ws.Range("N2").Value = -245.714285
# Row 7
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 200
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -88
$ws.Range("N7").ClearContents()
# Row 8
$ws.Range("H8").Value = 200
$ws.Range("I8").Value = 200
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 200
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -61
$ws.Range("N8").ClearContents()
# Row 62
$ws.Range("H62").Value = 30000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31372
$ws.Range("M62").ClearContents()
# Row 65
$ws.Range("H65").Value = 30000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96864
$ws.Range("M65").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 63
$ws.Range("H63").Value = 29999.5
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 29999.5
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 29999.5
$ws.Range("N63").Value = -31497.5
$ws.Range("M63").ClearContents()
# Row 66
$ws.Range("H66").Value = 29999.5
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 29999.5
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 89998.5
$ws.Range("N66").Value = -97486.5
$ws.Range("M66").ClearContents()
# Row 82
$ws.Range("H82").Value = 4386.625
$ws.Range("I82").Value = 4018.6
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 4018.6
$ws.Range("L82").Value = 5000
$ws.Range("M82").Value = -3657.6
$ws.Range("N82").Value = -5722
# Row 85
$ws.Range("H85").Value = 4386.625
$ws.Range("I85").Value = 4018.6
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 4018.6
$ws.Range("L85").Value = 5000
$ws.Range("M85").Value = -2770.6
$ws.Range("N85").Value = -7496

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 27061.75
$ws.Range("J62").Value = 51374
$ws.Range("L62").Value = 51374
$ws.Range("N62").Value = -52622
# Row 65
$ws.Range("H65").Value = 27061.75
$ws.Range("J65").Value = 51374
$ws.Range("L65").Value = 256870
$ws.Range("N65").Value = -263110
